$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2112211221122112
$ws.Range("C2").Value = 0.5445544554455446
$ws.Range("J2").Value = 0.009900990099009901
$ws.Range("P2").Value = 0.1452145214521452
$ws.Range("S2").Value = 0.0891089108910891
$ws.Range("B3").Value = 0.005813953488372093
$ws.Range("C3").Value = 0.05232558139534884
$ws.Range("J3").Value = 0.03488372093023256
$ws.Range("P3").Value = 0.7267441860465116
$ws.Range("S3").Value = 0.1802325581395349
$ws.Range("J4").Value = 0.01886792452830189
$ws.Range("P4").Value = 0.7547169811320755
$ws.Range("S4").Value = 0.2264150943396226
$ws.Range("B6").Value = 0.04779411764705882
$ws.Range("D6").Value = 0.007352941176470588
$ws.Range("F6").Value = 0.05882352941176471
$ws.Range("J6").Value = 0.2058823529411765
$ws.Range("O6").Value = 0.02205882352941177
$ws.Range("Q6").Value = 0.1544117647058824
$ws.Range("R6").Value = 0.05882352941176471
$ws.Range("S6").Value = 0.4448529411764706
$ws.Range("B7").Value = 0.08979591836734693
$ws.Range("D7").Value = 0.02448979591836735
$ws.Range("E7").Value = 0.004081632653061225
$ws.Range("F7").Value = 0.05714285714285714
$ws.Range("J7").Value = 0.08163265306122448
$ws.Range("O7").Value = 0.0326530612244898
$ws.Range("Q7").Value = 0.1877551020408163
$ws.Range("R7").Value = 0.0653061224489796
$ws.Range("S7").Value = 0.4571428571428571
$ws.Range("B8").Value = 0.09151414309484193
$ws.Range("D8").Value = 0.01830282861896839
$ws.Range("F8").Value = 0.05657237936772046
$ws.Range("J8").Value = 0.1064891846921797
$ws.Range("O8").Value = 0.01331114808652246
$ws.Range("Q8").Value = 0.1747088186356073
$ws.Range("R8").Value = 0.08153078202995008
$ws.Range("S8").Value = 0.4575707154742096
$ws.Range("B9").Value = 0.08118081180811808
$ws.Range("D9").Value = 0.007380073800738007
$ws.Range("F9").Value = 0.04797047970479705
$ws.Range("J9").Value = 0.07380073800738007
$ws.Range("O9").Value = 0.01845018450184502
$ws.Range("Q9").Value = 0.1955719557195572
$ws.Range("R9").Value = 0.07749077490774908
$ws.Range("S9").Value = 0.4981549815498155
$ws.Range("B10").Value = 0.08967596081386586
$ws.Range("D10").Value = 0.02637528259231349
$ws.Range("E10").Value = 0.001507159005275057
$ws.Range("F10").Value = 0.08590806330067823
$ws.Range("J10").Value = 0.1175584024114544
$ws.Range("O10").Value = 0.01205727204220045
$ws.Range("Q10").Value = 0.2019593067068576
$ws.Range("R10").Value = 0.06028636021100226
$ws.Range("S10").Value = 0.4046721929163527
$ws.Range("G11").Value = 0.1316614420062696
$ws.Range("J11").Value = 0.05642633228840126
$ws.Range("K11").Value = 0.1724137931034483
$ws.Range("L11").Value = 0.6144200626959248
$ws.Range("S11").Value = 0.02507836990595611
$ws.Range("G12").Value = 0.803921568627451
$ws.Range("J12").Value = 0.1323529411764706
$ws.Range("K12").Value = 0.004901960784313725
$ws.Range("L12").Value = 0.01470588235294118
$ws.Range("S12").Value = 0.04411764705882353
$ws.Range("G13").Value = 0.7272727272727273
$ws.Range("J13").Value = 0.1948051948051948
$ws.Range("S13").Value = 0.07792207792207792
$ws.Range("F15").Value = 0.01626016260162602
$ws.Range("H15").Value = 0.1869918699186992
$ws.Range("I15").Value = 0.0975609756097561
$ws.Range("J15").Value = 0.3089430894308943
$ws.Range("K15").Value = 0.04471544715447155
$ws.Range("M15").Value = 0.02845528455284553
$ws.Range("O15").Value = 0.04878048780487805
$ws.Range("S15").Value = 0.2682926829268293
$ws.Range("F16").Value = 0.02427184466019417
$ws.Range("H16").Value = 0.2233009708737864
$ws.Range("I16").Value = 0.0825242718446602
$ws.Range("J16").Value = 0.3203883495145631
$ws.Range("K16").Value = 0.06796116504854369
$ws.Range("M16").Value = 0.05825242718446602
$ws.Range("O16").Value = 0.07281553398058252
$ws.Range("S16").Value = 0.1504854368932039
$ws.Range("F17").Value = 0.02352941176470588
$ws.Range("H17").Value = 0.1686274509803922
$ws.Range("I17").Value = 0.1098039215686274
$ws.Range("J17").Value = 0.3745098039215686
$ws.Range("K17").Value = 0.08627450980392157
$ws.Range("M17").Value = 0.02352941176470588
$ws.Range("O17").Value = 0.07843137254901961
$ws.Range("S17").Value = 0.1352941176470588
$ws.Range("F18").Value = 0.01621621621621622
$ws.Range("H18").Value = 0.2594594594594595
$ws.Range("I18").Value = 0.05945945945945946
$ws.Range("J18").Value = 0.3405405405405406
$ws.Range("K18").Value = 0.0918918918918919
$ws.Range("M18").Value = 0.01621621621621622
$ws.Range("O18").Value = 0.08108108108108109
$ws.Range("S18").Value = 0.1351351351351351
$ws.Range("F19").Value = 0.01090909090909091
$ws.Range("H19").Value = 0.2266666666666667
$ws.Range("I19").Value = 0.1024242424242424
$ws.Range("J19").Value = 0.3454545454545455
$ws.Range("K19").Value = 0.1024242424242424
$ws.Range("M19").Value = 0.02606060606060606
$ws.Range("O19").Value = 0.06
$ws.Range("S19").Value = 0.1260606060606061
